$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final table of test results (rows 3-19), replacing the old "Serial Test" rows
# and inserting new rows so row 2 and row 20 stay unused (matching target dimension A1:D19).
$data = @(
    @(3,  "Poker Test",                         0.855,               "Passed", 17),
    @(4,  "Two Bit Test",                       0.385,               "Passed", 27),
    @(5,  "Gap Test",                           0.625,               "Passed", 29),
    @(6,  "Turning Point Test",                 0.227,               "Passed", 2),
    @(7,  "Autocorrelation Test",               0.575,               "Passed", 19),
    @(8,  "Hamming Weight Test",                0.9379999999999999,  "Passed", 7),
    @(9,  "Monobit",                            0.589,               "Passed", 0),
    @(10, "Frequency Within Block",             0.8070000000000001,  "Passed", 1),
    @(11, "Runs",                               0.253,               "Passed", 5),
    @(12, "Longest Run Ones In A Block",        0.585,               "Passed", 19),
    @(13, "Discrete Fourier Transform",         0.783,               "Passed", 1),
    @(14, "Non Overlapping Template Matching",  1,                   "Passed", 50),
    @(15, "Serial Test",                        0.915,               "Passed", 438),
    @(16, "Approximate Entropy",                0.741,               "Passed", 531),
    @(17, "Cumulative Sums",                    0.739,               "Passed", 20),
    @(18, "Random Excursion",                   0.038,               "Failed", 57),
    @(19, "Random Excursion Variant",           0.432,               "Passed", 1)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

# Remove the old row 21 (now unused since the table only spans rows 1-19)
$ws.Rows.Item(21).Delete()
